$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header M1 from "Uraian Ruang" to "Kode Gedung" (add scan qrcode page related column)
$ws.Range("M1").Value = "Kode Gedung"

# Update the active selection to K8 as recorded in the saved workbook
$ws.Range("K8").Select()
